# Auto-generated script to apply market-data refresh values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 302.66666
$ws.Range("I2").Value = 54
$ws.Range("K2").Value = 54
$ws.Range("M2").Value = 59
$ws.Range("H9").Value = 238.28572
$ws.Range("I9").Value = 88.8
$ws.Range("K9").Value = 88.8
$ws.Range("M9").Value = 80.2
$ws.Range("H43").Value = 2227.375
$ws.Range("I43").Value = 1231.6666
$ws.Range("K43").Value = 1231.6666
$ws.Range("M43").Value = -1162.6666
$ws.Range("H69").Value = 6763.5
$ws.Range("J69").Value = 6763.5
$ws.Range("L69").Value = 20290.5
$ws.Range("N69").Value = -22038.5
$ws.Range("H72").Value = 6763.5
$ws.Range("J72").Value = 6763.5
$ws.Range("L72").Value = 60871.5
$ws.Range("N72").Value = -69607.5
$ws.Range("H86").Value = 4690
$ws.Range("I86").Value = 4598.75
$ws.Range("K86").Value = 4598.75
$ws.Range("M86").Value = -3475.75
$ws.Range("H89").Value = 4690
$ws.Range("I89").Value = 4598.75
$ws.Range("K89").Value = 22993.75
$ws.Range("M89").Value = -17377.75
$ws.Range("H106").Value = 2921.0908
$ws.Range("I106").Value = 2921.0908
$ws.Range("K106").Value = 2921.0908
$ws.Range("M106").Value = -2290.0908
$ws.Range("H116").Value = 28317.727
$ws.Range("I116").Value = 7427.5
$ws.Range("K116").Value = 7427.5
$ws.Range("M116").Value = -3985.5
$ws.Range("H138").Value = 3762.5386
$ws.Range("I138").Value = 6043
$ws.Range("J138").Value = 3572.5
$ws.Range("K138").Value = 18129
$ws.Range("L138").Value = 10717.5
$ws.Range("M138").Value = -12989
$ws.Range("N138").Value = -20997.5
$ws.Range("H141").Value = 8607
$ws.Range("J141").Value = 12994.833
$ws.Range("L141").Value = 38984.499
$ws.Range("N141").Value = -49344.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 370.4
$ws.Range("I5").Value = 100.71429
$ws.Range("J5").Value = 999.6667
$ws.Range("K5").Value = 100.71429
$ws.Range("L5").Value = 999.6667
$ws.Range("M5").Value = 11.28570999999999
$ws.Range("N5").Value = -1223.6667
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H74").Value = 3277184.5
$ws.Range("I74").Value = 4276587.5
$ws.Range("K74").Value = 4276587.5
$ws.Range("M74").Value = -4275713.5
$ws.Range("H77").Value = 3277184.5
$ws.Range("I77").Value = 4276587.5
$ws.Range("K77").Value = 21382937.5
$ws.Range("M77").Value = -21378569.5
$ws.Range("H133").Value = 86418.836
$ws.Range("J133").Value = 86418.836
$ws.Range("L133").Value = 86418.836
$ws.Range("N133").Value = -91478.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 370.4
$ws.Range("I4").Value = 100.71429
$ws.Range("J4").Value = 999.6667
$ws.Range("K4").Value = 100.71429
$ws.Range("L4").Value = 999.6667
$ws.Range("M4").Value = 14.28570999999999
$ws.Range("N4").Value = -1229.6667
$ws.Range("H86").Value = 1293.2
$ws.Range("I86").Value = 703
$ws.Range("J86").Value = 1686.6666
$ws.Range("K86").Value = 703
$ws.Range("L86").Value = 1686.6666
$ws.Range("M86").Value = 420
$ws.Range("N86").Value = -3932.6666
$ws.Range("H89").Value = 1293.2
$ws.Range("I89").Value = 703
$ws.Range("J89").Value = 1686.6666
$ws.Range("K89").Value = 3515
$ws.Range("L89").Value = 8433.333000000001
$ws.Range("M89").Value = 2101
$ws.Range("N89").Value = -19665.333
$ws.Range("H94").Value = 2457.6667
$ws.Range("J94").Value = 1502.25
$ws.Range("L94").Value = 1502.25
$ws.Range("N94").Value = -2404.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2681.6924
$ws.Range("I62").Value = 2733
$ws.Range("J62").Value = 2399.5
$ws.Range("K62").Value = 2733
$ws.Range("L62").Value = 2399.5
$ws.Range("M62").Value = -2109
$ws.Range("N62").Value = -3647.5
$ws.Range("H65").Value = 2681.6924
$ws.Range("I65").Value = 2733
$ws.Range("J65").Value = 2399.5
$ws.Range("K65").Value = 13665
$ws.Range("L65").Value = 11997.5
$ws.Range("M65").Value = -10545
$ws.Range("N65").Value = -18237.5
$ws.Range("H92").Value = 86666.664
$ws.Range("J92").Value = 86666.664
$ws.Range("L92").Value = 86666.664
$ws.Range("N92").Value = -91658.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1302.3
$ws.Range("I113").Value = 825.7143
$ws.Range("K113").Value = 2477.1429
$ws.Range("M113").Value = -307.1428999999998
$ws.Range("H137").Value = 5897.3794
$ws.Range("I137").Value = 1848.0769
$ws.Range("J137").Value = 9187.4375
$ws.Range("K137").Value = 5544.2307
$ws.Range("L137").Value = 27562.3125
$ws.Range("M137").Value = -444.2307000000001
$ws.Range("N137").Value = -37762.3125
$ws.Range("H140").Value = 2295.8462
$ws.Range("I140").Value = 1734.5
$ws.Range("K140").Value = 5203.5
$ws.Range("M140").Value = -23.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7887.2
$ws.Range("I80").Value = 6034.4
$ws.Range("K80").Value = 6034.4
$ws.Range("M80").Value = -5036.4
$ws.Range("H83").Value = 7887.2
$ws.Range("I83").Value = 6034.4
$ws.Range("K83").Value = 30172
$ws.Range("M83").Value = -25180
$ws.Range("H95").Value = 28398.6
$ws.Range("J95").Value = 28398.6
$ws.Range("L95").Value = 28398.6
$ws.Range("N95").Value = -33890.6
$ws.Range("H122").Value = 70175.88
$ws.Range("J122").Value = 10374.75
$ws.Range("L122").Value = 31124.25
$ws.Range("N122").Value = -36024.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1812.625
$ws.Range("I68").Value = 1795.4286
$ws.Range("K68").Value = 1795.4286
$ws.Range("M68").Value = -1046.4286
$ws.Range("H71").Value = 1812.625
$ws.Range("I71").Value = 1795.4286
$ws.Range("K71").Value = 8977.143
$ws.Range("M71").Value = -5233.143
$ws.Range("H82").Value = 2624.5454
$ws.Range("I82").Value = 3052.8572
$ws.Range("J82").Value = 1875
$ws.Range("K82").Value = 3052.8572
$ws.Range("L82").Value = 1875
$ws.Range("M82").Value = -2691.8572
$ws.Range("N82").Value = -2597
$ws.Range("H85").Value = 2624.5454
$ws.Range("I85").Value = 3052.8572
$ws.Range("J85").Value = 1875
$ws.Range("K85").Value = 3052.8572
$ws.Range("L85").Value = 1875
$ws.Range("M85").Value = -1804.8572
$ws.Range("N85").Value = -4371
$ws.Range("H133").Value = 77732.5
$ws.Range("J133").Value = 70396
$ws.Range("L133").Value = 70396
$ws.Range("N133").Value = -75456

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2008.3334
$ws.Range("I6").Value = 125
$ws.Range("J6").Value = 2950
$ws.Range("K6").Value = 125
$ws.Range("L6").Value = 2950
$ws.Range("M6").Value = -10
$ws.Range("N6").Value = -3180
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = ""
$ws.Range("H62").Value = 10052.111
$ws.Range("J62").Value = 15292.25
$ws.Range("L62").Value = 15292.25
$ws.Range("N62").Value = -16540.25
$ws.Range("H65").Value = 10052.111
$ws.Range("J65").Value = 15292.25
$ws.Range("L65").Value = 76461.25
$ws.Range("N65").Value = -82701.25
$ws.Range("H96").Value = 1624.875
$ws.Range("I96").Value = 1583.1666
$ws.Range("J96").Value = 1750
$ws.Range("K96").Value = 1583.1666
$ws.Range("L96").Value = 1750
$ws.Range("M96").Value = -210.1666
$ws.Range("N96").Value = -4496
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H122").Value = 70906.3
$ws.Range("I122").Value = 5924.5835
$ws.Range("K122").Value = 17773.7505
$ws.Range("M122").Value = -15323.7505
$ws.Range("H136").Value = 11068542
$ws.Range("I136").Value = 2718728.2
$ws.Range("K136").Value = 8156184.600000001
$ws.Range("M136").Value = -8153634.600000001
